$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.956.08'
$ws.Range("E2").Value = '  +0.89%  '
$ws.Range("D3").Value = '2.420.00'
$ws.Range("E3").Value = '  +0.78%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '552.37'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.12%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.49'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.76%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.588'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.19%  '
$ws.Range("E9").Value = '  -1.53%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '5.69'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.78%  '
$ws.Range("E11").Value = '  -1.55%  '
$ws.Range("E12").Value = '  -1.02%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '25.52'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.71%  '
$ws.Range("D14").Value = '2.849.88'
$ws.Range("E14").Value = '  +0.75%  '
$ws.Range("D15").Value = '59.876.94'
$ws.Range("E15").Value = '  +0.87%  '
$ws.Range("E16").Value = '  -0.51%  '
$ws.Range("D17").Value = '2.430.16'
$ws.Range("E17").Value = '  +0.50%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.38'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.16%  '
$ws.Range("E19").Value = '  -0.15%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '329.69'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.37%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.65'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.23%  '
$ws.Range("E22").Value = '  +0.07%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '66.51'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.86%  '
$ws.Range("E24").Value = '  +2.04%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.64'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.79%  '
$ws.Range("E26").Value = '  +0.03%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.38'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.34%  '
$ws.Range("E28").Value = '  +0.87%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '169.25'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.11%  '
$ws.Range("E31").Value = '  -1.94%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '18.67'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.11%  '
$ws.Range("E33").Value = '  -0.77%  '
$ws.Range("E35").Value = '  +1.10%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.24'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.95%  '
$ws.Range("E37").Value = '  +0.16%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.61'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.63%  '
$ws.Range("E39").Value = '  -2.52%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '314.17'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.06%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.68'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.66%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '138.61'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.67%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0967'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.76%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0520'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.72%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '19.59'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.26%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.580'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.03%  '
$ws.Range("E47").Value = '  -0.37%  '
$ws.Range("E48").Value = '  -5.19%  '
$ws.Range("E49").Value = '  +0.32%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '11.07'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.34%  '
$ws.Range("E51").Value = '  +0.21%  '
